$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# IOT sheet: update marks so that the attendance/CP line (row 6) gets a score
# and the overall totals on row 12 roll up rows 6-11 (instead of 7-10).
# ---------------------------------------------------------------------------
$iot = $wb.Worksheets.Item("IOT")
$iot.Range("C6").Value = 20
$iot.Range("D6").Value = 20
$iot.Range("G6").Value = 20

$iot.Range("C11").Style = $iot.Range("C9").Style

$iot.Range("C12").Formula = "=SUM(C6:C11)"
$iot.Range("E12").Formula = "=ROUND(SUM(E6:E11),2)"

# ---------------------------------------------------------------------------
# TDD sheet: move the submitted CP mark from S5 (row 11) to S4 (row 10) and
# bump up the total-marks denominator for that line.
# ---------------------------------------------------------------------------
$tdd = $wb.Worksheets.Item("TDD")
$tdd.Range("D10").Value = 100
$tdd.Range("G10").Value = 81
$tdd.Range("G11").ClearContents()

# ---------------------------------------------------------------------------
# OR sheet: shared SUM formula for F10 / F11 (cosmetic formula metadata only
# -- re-enter the same formula so the dependency stays correct).
# ---------------------------------------------------------------------------
$or = $wb.Worksheets.Item("OR")
$or.Range("F10").Formula = "=SUM(G10:M10)"

# ---------------------------------------------------------------------------
# WE sheet: re-apply the header style on the merged "Marks" banner (G4:J4).
# ---------------------------------------------------------------------------
$we = $wb.Worksheets.Item("WE")
$weHeader = $we.Range("G4:J4")
$weHeader.Font.Name = "Times New Roman"
$weHeader.Font.Size = 12
$weHeader.Font.Bold = $true
$weHeader.HorizontalAlignment = -4108
$weHeader.VerticalAlignment = -4108
$weHeader.WrapText = $true

# ---------------------------------------------------------------------------
# GPA sheet: credit-hours row (9) becomes a live formula driven off row 8
# (previously hard-coded 3's), which ripples down into the weighted totals.
# ---------------------------------------------------------------------------
$gpa = $wb.Worksheets.Item("GPA")
$gpa.Range("D9").Formula = "=IF(D8=0, 0,3)"
$gpa.Range("E9").Formula = "=IF(E8=0, 0,3)"
$gpa.Range("F9").Formula = "=IF(F8=0, 0,3)"
$gpa.Range("G9").Formula = "=IF(G8=0, 0,3)"
$gpa.Range("H9").Formula = "=IF(H8=0, 0,3)"
$gpa.Range("I9").Formula = "=IF(I8=0, 0,3)"

# ---------------------------------------------------------------------------
# Restore on-screen selections to match where everyone ended up working.
# The last sheet selected/activated becomes the active tab, so TDD goes last.
# ---------------------------------------------------------------------------
$or.Range("G10").Select() | Out-Null

$we.Range("G8:G10").Select() | Out-Null

$hci = $wb.Worksheets.Item("HCI")
$hci.Range("G9").Select() | Out-Null

$iot.Range("G6").Select() | Out-Null

$gpa.Range("J9").Select() | Out-Null

$tdd.Range("G7").Select() | Out-Null
$tdd.Activate()
